# Apply the two edits described by the diff:
# 1. Remove the whole first paragraph ("I will report what I have done ...
#    ... for discussion.") including its paragraph mark.
# 2. Merge the trailing "achieved" run (and the space run before it) into
#    the preceding run so the sentence reads as a single run of text.

$d = $word.ActiveDocument

# --- Edit 1: delete the entire first paragraph ------------------------
$first = $d.Paragraphs.First
$first.Range.Delete()

# --- Edit 2: collapse the trailing " " run and the "achieved" run into
#     the preceding " ... what benefit can be" run. Searching for just
#     " achieved" (the text spanning those two small runs) keeps the
#     unrelated " of sensor data" run untouched while merging the rest. --
$d.Content.Find.Execute(
    " achieved",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    " achieved",
    2)
